# Daily attendance processing - 2025-12-03 12:45:24
# Updates "Recorded By" attendee orderings, a couple of recorded attendance
# counts, and two average-attendance percentage figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Session Analysis Results")

# --- "Recorded By" column: attendee lists re-ordered (same people, new order) ---
$ws.Range("G2").Value  = "servinaz@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg, System, Amira.Sobhy@med.asu.edu.eg, gehanadel@med.asu.edu.eg"
$ws.Range("G3").Value  = "hend_mahmoud@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg, majorelle.magdy@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg, System"
$ws.Range("G4").Value  = "hend_mahmoud@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg, servinaz@med.asu.edu.eg, majorelle.magdy@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg, gehanadel@med.asu.edu.eg"
$ws.Range("G5").Value  = "eman.tantawi@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg"
$ws.Range("G6").Value  = "majorelle.magdy@med.asu.edu.eg, alshimaa.atef@med.asu.edu.egm, Mohammedeltanany@med.asu.edu.eg"
$ws.Range("G7").Value  = "NadaMohamed@med.asu.edu.eg, Amera.a.saad@med.asu.edu.eg, menna-alah.mohamed@asu.edu.eg, AbeerRagheb@med.asu.edu.eg, Kerelos.zareef@med.asu.edu.eg, Fatmaelhady@med.asu.edu.eg, lamiaa.ossama@med.asu.edu.eg"
$ws.Range("G12").Value = "Madeha.Saeed@med.asu.edu.eg, amira.m.ibrahim@med.asu.edu.eg, dina.adel@med.asu.edu.eg, yassmina.fattoh@med.asu.edu.eg, Eman.m.abosakaya@med.asu.edu.eg, Marina.youhana@med.asu.edu.eg"
$ws.Range("G28").Value = "maryam.ashraf@med.asu.edu.eg, Aya_hamed@med.asu.edu.eg"
$ws.Range("G30").Value = "yassmen.ahmed@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg, shorokmohamed@med.asu.edu.eg, aya.hanafy@med.asu.edu.eg"

# --- Updated recorded-attendance counts ("x/251" stays plain text) ---
$ws.Range("H14").Value = "97/251"
$ws.Range("H21").Value = "66/251"

# --- Updated average-attendance percentages ---
# NOTE: assigning a literal "26.6%" string straight to .Value gets
# auto-recognized by Excel as a number and reformatted into a numeric
# percentage cell, which would change the cell's type/style. To keep it as
# literal text (matching the rest of the report, which stores percentages as
# text), build the string with a formula so it is never number-parsed, then
# flatten the formula down to a static value via copy / paste-values. This
# preserves the original cell style and text type.
$cL10 = $ws.Range("L10")
$cL10.Formula = "=""26.6""&""%"""
$cL10.Copy()
$cL10.PasteSpecial(-4163)

$cS15 = $ws.Range("S15")
$cS15.Formula = "=""26.6""&""%"""
$cS15.Copy()
$cS15.PasteSpecial(-4163)

$excel.CutCopyMode = 0
